# "Generate Report for Handoff"
#
# The localization CI run moved from "In Translation" to "Ready for
# handoff": the status text + the two xliff-generation timestamps change,
# and the Status/"Latest HO Xliff Generate Date"/"Latest Handoff Datetime"
# columns get a little wider to fit the new (longer) status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # Status column
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # Status column

# --- Handoff timestamps (new xliff generated a little later) ---------------
$wsOverview.Range("G2").Value = "2016-08-22 22:38:51" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-22 22:38:51" # Latest Handoff Datetime
$wsZhCn.Range("H2").Value     = "2016-08-22 22:38:47" # Latest Handoff Datetime

# --- Widen the Status / date columns to fit "Ready for handoff" ------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33   # column C
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33   # column C
